$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DPbES")

# --- New rows: crude oil (15), heavy or residual fuel oil (16), municipal solid waste (17) ---
# (added first so the new shared-string entries land in the same order as the
# authored workbook: crude oil, heavy or residual fuel oil, municipal solid
# waste, then the new column-A header text)
$ws.Range("A15").Value = "crude oil"
$ws.Range("B15").Formula = "=B11"
$ws.Range("C15:AK15").Formula = "=C11"

$ws.Range("A16").Value = "heavy or residual fuel oil"
$ws.Range("B16").Formula = "=B11"
$ws.Range("C16:AK16").Formula = "=C11"

$ws.Range("A17").Value = "municipal solid waste"
$ws.Range("B17").Formula = "=B9"
$ws.Range("C17:AK17").Formula = "=C9"

# --- Header cell A1: new title "Dispatch Priority (dimensionless)" ---
$ws.Range("A1").Value = "Dispatch Priority (dimensionless)"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 30

# --- Column A slightly wider to fit new content ---
$ws.Columns.Item(1).ColumnWidth = 23
